$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.478.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.17%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.958.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'542.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -4.74%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'151.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -6.85%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.14%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.80%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.965.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -5.88%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "'  -3.36%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.476.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -6.07%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.87%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'61.557.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.18%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'23.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.45%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.963.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -6.01%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -4.85%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.96%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'381.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.28%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -5.75%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.18%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.49%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.60%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.083.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -6.57%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.89%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.01%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0931"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -7.82%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'8.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.01%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -4.45%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'20.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.25%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'159.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.14%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.39%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -5.02%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.29%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.48%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -6.69%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.408.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -9.57%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.16%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'22.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.48%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.662"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.45%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0592"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.56%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.29%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -3.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'4.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -8.31%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0957"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.02%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Bittensor"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'267.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -6.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'InjectiveProtocol"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'19.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.58%  "
$ws.Range("E51").Style = "Normal"
